$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.412774443626404
$ws.Range("B1").Value = 2.20635199546814
$ws.Range("C1").Value = 4.958789348602295
$ws.Range("D1").Value = 3.319113492965698
$ws.Range("E1").Value = 1.186159014701843
